$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# New "id" column header in Sheet2!A1 (adds a new shared string "id").
$ws2.Range("A1").Value = "id"

# Move the active selection on Sheet2 from C6 to A2.
$ws2.Activate()
[void]$ws2.Range("A2").Select()
